# Fixing name of Sectors to be aligned with Baseline
#
# 1) Rename the four element/commodity header labels (row 3, columns D:G)
#    on every worksheet from their abbreviations to full names:
#      Nd -> Neodymium
#      Dy -> Dysprosium
#      Cu -> Copper ores and concentrates
#      Si -> Raw silicon
#
# 2) A small number of worksheets' G7 cell carries a value that shifted by
#    a single ULP (floating point recalculation noise) between the
#    baseline and this commit. Re-apply those exact values.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Cells.Item(3, 4).Value = "Neodymium"
    $ws.Cells.Item(3, 5).Value = "Dysprosium"
    $ws.Cells.Item(3, 6).Value = "Copper ores and concentrates"
    $ws.Cells.Item(3, 7).Value = "Raw silicon"
}

$g7Fixes = @{
    8   = -3996.147468921459
    11  = -48234.30674689904
    25  = -78263307.30023907
    30  = -155685484.1210942
    36  = -256209997.4367262
    40  = -405460198.430169
    43  = -739082373.7167541
    44  = -929058414.7353605
    53  = -2455549392.221612
    55  = -2304601632.933053
    62  = -718953063.5615468
    63  = -576730941.9844749
    75  = -109864769.3578359
    80  = -193429155.7372571
    95  = -285754767.840836
    100 = -290258297.5629358
}

foreach ($idx in $g7Fixes.Keys) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Cells.Item(7, 7).Value = $g7Fixes[$idx]
}
